# Update the "mines - ..." build/version tag throughout the workbook.
#
# Old:  mines - January 30 (built on February 02 2026 12.49.33 EST)
# New:  mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)

$wb = $excel.ActiveWorkbook

$oldTag = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newTag = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet -------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: " + $newTag

$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Blackwater Coal Mine, Australia, M0011, version '" + $newTag + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet -----------------------------
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 23; $row++) {
    $cell = $data.Cells.Item($row, 19)  # column S = build_version
    if ($cell.Value() -eq $oldTag) {
        $cell.Value = $newTag
    }
}
